$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60797
$ws.Range("B3").Value = 131419
$ws.Range("B4").Value = 123272
$ws.Range("B5").Value = 13340
$ws.Range("B6").Value = 27529
$ws.Range("B8").Value = 3071496
$ws.Range("B9").Value = 217900
$ws.Range("B10").Value = 29884
$ws.Range("B11").Value = 626239
$ws.Range("B12").Value = 323841
$ws.Range("B14").Value = 183330
$ws.Range("B15").Value = 767338
$ws.Range("B16").Value = 3915
$ws.Range("B17").Value = 363732
$ws.Range("B18").Value = 1003746
$ws.Range("B19").Value = 12680
$ws.Range("B21").Value = 1144
$ws.Range("B22").Value = 310572
$ws.Range("B23").Value = 199882
$ws.Range("B25").Value = 14930183
$ws.Range("B26").Value = 228
$ws.Range("B27").Value = 407827
$ws.Range("B28").Value = 13345
$ws.Range("B29").Value = 142874
$ws.Range("B30").Value = 4069
$ws.Range("B31").Value = 25159
$ws.Range("B32").Value = 16971
$ws.Range("B34").Value = 1266249
$ws.Range("B35").Value = 6521
$ws.Range("B36").Value = 4847
$ws.Range("B37").Value = 1222949
$ws.Range("B38").Value = 102571
$ws.Range("B39").Value = 2934611
$ws.Range("B40").Value = 3847
$ws.Range("B42").Value = 30130
$ws.Range("B43").Value = 260535
$ws.Range("B44").Value = 46231
$ws.Range("B45").Value = 339412
$ws.Range("B46").Value = 111654
$ws.Range("B47").Value = 67982
$ws.Range("B48").Value = 1639248
$ws.Range("B49").Value = 255181
$ws.Range("B50").Value = 11268
$ws.Range("B52").Value = 268561
$ws.Range("B53").Value = 393048
$ws.Range("B54").Value = 232905
$ws.Range("B55").Value = 69997
$ws.Range("B57").Value = 3707
$ws.Range("B58").Value = 123781
$ws.Range("B59").Value = 18467
$ws.Range("B60").Value = 260139
$ws.Range("B61").Value = 125
$ws.Range("B62").Value = 87798
$ws.Range("B63").Value = 5767541
$ws.Range("B65").Value = 5914
$ws.Range("B66").Value = 315913
$ws.Range("B67").Value = 3484755
$ws.Range("B68").Value = 92828
$ws.Range("B69").Value = 352027
$ws.Range("B71").Value = 231289
$ws.Range("B72").Value = 22468
$ws.Range("B73").Value = 3738
$ws.Range("B74").Value = 13829
$ws.Range("B75").Value = 13149
$ws.Range("B76").Value = 216964
$ws.Range("B77").Value = 785967
$ws.Range("B78").Value = 6497
$ws.Range("B79").Value = 21077410
$ws.Range("B80").Value = 1691658
$ws.Range("B81").Value = 2591609
$ws.Range("B82").Value = 1091954
$ws.Range("B83").Value = 251087
$ws.Range("B84").Value = 838767
$ws.Range("B85").Value = 4070400
$ws.Range("B86").Value = 46194
$ws.Range("B87").Value = 617890
$ws.Range("B88").Value = 716923
$ws.Range("B89").Value = 388111
$ws.Range("B90").Value = 161393
$ws.Range("B91").Value = 125519
$ws.Range("B93").Value = 280536
$ws.Range("B94").Value = 96958
$ws.Range("B95").Value = 1072
$ws.Range("B96").Value = 120736
$ws.Range("B97").Value = 530217
$ws.Range("B98").Value = 10749
$ws.Range("B99").Value = 2113
$ws.Range("B100").Value = 178927
$ws.Range("B101").Value = 2956
$ws.Range("B102").Value = 252699
$ws.Range("B103").Value = 67850
$ws.Range("B104").Value = 38116
$ws.Range("B105").Value = 34143
$ws.Range("B106").Value = 424376
$ws.Range("B107").Value = 32665
$ws.Range("B108").Value = 13998
$ws.Range("B109").Value = 30411
$ws.Range("B111").Value = 18542
$ws.Range("B112").Value = 1216
$ws.Range("B113").Value = 2355985
$ws.Range("B114").Value = 251820
$ws.Range("B115").Value = 2471
$ws.Range("B116").Value = 41524
$ws.Range("B117").Value = 97930
$ws.Range("B118").Value = 512656
$ws.Range("B119").Value = 70052
$ws.Range("B120").Value = 49323
$ws.Range("B121").Value = 359610
$ws.Range("B122").Value = 1558549
$ws.Range("B123").Value = 2633
$ws.Range("B124").Value = 6989
$ws.Range("B125").Value = 5286
$ws.Range("B126").Value = 165273
$ws.Range("B127").Value = 114905
$ws.Range("B128").Value = 198572
$ws.Range("B129").Value = 841636
$ws.Range("B130").Value = 365975
$ws.Range("B131").Value = 11273
$ws.Range("B132").Value = 288974
$ws.Range("B133").Value = 1818689
$ws.Range("B134").Value = 1073555
$ws.Range("B135").Value = 2811951
$ws.Range("B136").Value = 838102
$ws.Range("B137").Value = 208877
$ws.Range("B138").Value = 1060895
$ws.Range("B139").Value = 4792354
$ws.Range("B140").Value = 25421
$ws.Range("B142").Value = 4585
$ws.Range("B146").Value = 2314
$ws.Range("B147").Value = 422316
$ws.Range("B148").Value = 40544
$ws.Range("B149").Value = 695875
$ws.Range("B151").Value = 4068
$ws.Range("B152").Value = 61268
$ws.Range("B153").Value = 384317
$ws.Range("B154").Value = 243719
$ws.Range("B156").Value = 14121
$ws.Range("B157").Value = 1588221
$ws.Range("B158").Value = 10613
$ws.Range("B159").Value = 3551262
$ws.Range("B160").Value = 117529
$ws.Range("B161").Value = 34082
$ws.Range("B162").Value = 10696
$ws.Range("B163").Value = 995595
$ws.Range("B164").Value = 667380
$ws.Range("B165").Value = 23121
$ws.Range("B166").Value = 1160
$ws.Range("B169").Value = 74900
$ws.Range("B170").Value = 2786
$ws.Range("B171").Value = 13068
$ws.Range("B172").Value = 12105
$ws.Range("B173").Value = 315600
$ws.Range("B174").Value = 4955594
$ws.Range("B175").Value = 32557444
$ws.Range("B176").Value = 42102
$ws.Range("B177").Value = 2146121
$ws.Range("B178").Value = 529220
$ws.Range("B179").Value = 4441644
$ws.Range("B180").Value = 209867
$ws.Range("B181").Value = 92724
$ws.Range("B183").Value = 202578
$ws.Range("B184").Value = 3030
$ws.Range("B185").Value = 299736
$ws.Range("B186").Value = 6414
$ws.Range("B187").Value = 91849
$ws.Range("B188").Value = 38357
